$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert a new banner row at the top (pushes existing rows down by one) ---
[void]$ws.Rows.Item(1).Insert()

# New banner text (note about MSK null flavor) in A1, merged across A1:D1
$noteText = "Note: MSK is a null flavor that means ""masked."" MSK is used when there is information available for the value, but it has not been provided to the sender due to security, privacy, or other reasons."
$ws.Range("A1").Value2 = $noteText
[void]$ws.Range("A1:D1").Merge()
$ws.Range("A1:D1").HorizontalAlignment = -4108   # xlCenter

# --- Remove the two stray empty trailing rows (formerly 87/88, now 88/89) ---
[void]$ws.Rows.Item(88).Delete()
[void]$ws.Rows.Item(88).Delete()

# --- Column width adjustments (column B and column D), bestFit is dropped automatically ---
$ws.Columns.Item(2).ColumnWidth = 33.4986979166666
$ws.Columns.Item(4).ColumnWidth = 78.1666666666667

# --- Update the view: drop the scrolled topLeftCell and move the active selection ---
[void]$ws.Range("E4").Select()

Write-Output "done"
